$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.697.16'
$ws.Range("E2").Value = '  -2.37%  '

$ws.Range("D3").Value = '1.888.74'
$ws.Range("E3").Value = '  -3.43%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.98'
$ws.Range("E5").Value = '  -1.89%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.16%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4565'
$ws.Range("E7").Value = '  -1.55%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3789'
$ws.Range("E8").Value = '  -3.84%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.47'
$ws.Range("E9").Value = '  -1.58%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07676'
$ws.Range("E10").Value = '  -2.69%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9578'
$ws.Range("E11").Value = '  -4.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.89'
$ws.Range("E12").Value = '  -2.14%  '

$ws.Range("D13").Value = '1.887.10'
$ws.Range("E13").Value = '  -3.75%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.941'
$ws.Range("E14").Value = '  -3.10%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.634'
$ws.Range("E15").Value = '  -3.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07016'
$ws.Range("E16").Value = '  -1.47%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  -0.19%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '82.43'
$ws.Range("E18").Value = '  -6.78%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000009468'
$ws.Range("E19").Value = '  -4.72%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.58'
$ws.Range("E20").Value = '  -3.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("D22").Value = '28.684.00'
$ws.Range("E22").Value = '  -2.63%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.325'
$ws.Range("E23").Value = '  -3.56%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.83'
$ws.Range("E24").Value = '  -3.53%  '

$ws.Range("D25").Value = '2.125.52'
$ws.Range("E25").Value = '  -3.33%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.088'
$ws.Range("E26").Value = '  -1.46%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '155.02'
$ws.Range("E27").Value = '  -2.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.89'
$ws.Range("E28").Value = '  -3.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.590'
$ws.Range("E29").Value = '  -6.60%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.46'
$ws.Range("E30").Value = '  -2.73%  '

$ws.Range("E31").Value = '  -3.75%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09217'

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8404'
$ws.Range("E33").Value = '  -5.91%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.039'
$ws.Range("E34").Value = '  -3.89%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.240'
$ws.Range("E35").Value = '  -7.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.055'
$ws.Range("E36").Value = '  -3.79%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05634'
$ws.Range("E37").Value = '  -3.11%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.143'
$ws.Range("E38").Value = '  -3.08%  '

$ws.Range("E39").Value = '  -0.12%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02022'
$ws.Range("E40").Value = '  -4.55%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.434'
$ws.Range("E41").Value = '  -5.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5460'
$ws.Range("E42").Value = '  -5.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1742'
$ws.Range("E43").Value = '  -4.44%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.000002917'
$ws.Range("E44").Value = '  -22.97%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.143'
$ws.Range("E45").Value = '  -6.71%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.696'
$ws.Range("E46").Value = '  +3.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5147'
$ws.Range("E47").Value = '  -4.14%  '

$ws.Range("E48").Value = '  -7.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.076'
$ws.Range("E49").Value = '  -5.40%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06760'
$ws.Range("E50").Value = '  -2.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.47'
$ws.Range("E51").Value = '  -2.96%  '
